# Apply updated crypto price/volume data to Sheet1 (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.361.38"
$ws.Range("E2").Value = "  -2.78%  "

$ws.Range("D3").Value = "2.417.25"
$ws.Range("E3").Value = "  -3.76%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'" + "512.48"
$ws.Range("E5").Value = "  -3.94%  "

$ws.Range("D6").Value = "'" + "131.16"
$ws.Range("E6").Value = "  -3.43%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -3.17%  "

$ws.Range("D9").Value = "2.420.28"
$ws.Range("E9").Value = "  -3.82%  "

$ws.Range("D10").Value = "'" + "0.0954"
$ws.Range("E10").Value = "  -6.00%  "

$ws.Range("E11").Value = "  -1.53%  "

$ws.Range("D12").Value = "'" + "5.19"
$ws.Range("E12").Value = "  -4.28%  "

$ws.Range("E13").Value = "  -4.42%  "

$ws.Range("D14").Value = "2.848.44"
$ws.Range("E14").Value = "  -3.66%  "

$ws.Range("D15").Value = "57.315.12"
$ws.Range("E15").Value = "  -2.71%  "

$ws.Range("D16").Value = "'" + "21.54"
$ws.Range("E16").Value = "  -5.60%  "

$ws.Range("E17").Value = "  -4.51%  "

$ws.Range("D18").Value = "2.427.25"
$ws.Range("E18").Value = "  -3.30%  "

$ws.Range("E19").Value = "  -6.24%  "

$ws.Range("D20").Value = "'" + "313.61"
$ws.Range("E20").Value = "  -3.14%  "

$ws.Range("E21").Value = "  -4.01%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").Value = "'" + "5.61"
$ws.Range("E23").Value = "  -5.18%  "

$ws.Range("D24").Value = "'" + "63.71"
$ws.Range("E24").Value = "  -2.20%  "

$ws.Range("E25").Value = "  -4.42%  "

$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("E27").Value = "  -2.97%  "

$ws.Range("E28").Value = "  -5.06%  "

$ws.Range("D29").Value = "'" + "169.31"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "0.0₃0721"
$ws.Range("E30").Value = "  -5.66%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'" + "1.66"
$ws.Range("E31").Value = "  -5.06%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'" + "6.19"
$ws.Range("E32").Value = "  -5.18%  "

$ws.Range("E33").Value = "  +1.21%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").Value = "'" + "0.998"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").Value = "'" + "17.68"
$ws.Range("E36").Value = "  -3.96%  "

$ws.Range("E37").Value = "  -7.29%  "

$ws.Range("E38").Value = "  -4.84%  "

$ws.Range("D39").Value = "'" + "36.21"
$ws.Range("E39").Value = "  -2.41%  "

$ws.Range("E40").Value = "  -5.91%  "

$ws.Range("D41").Value = "'" + "0.770"
$ws.Range("E41").Value = "  -3.89%  "

$ws.Range("E42").Value = "  -6.47%  "

$ws.Range("D43").Value = "'" + "265.00"
$ws.Range("E43").Value = "  -6.16%  "

$ws.Range("D44").Value = "'" + "4.89"
$ws.Range("E44").Value = "  -2.86%  "

$ws.Range("E45").Value = "  -4.07%  "

$ws.Range("D46").Value = "'" + "121.52"
$ws.Range("E46").Value = "  -6.27%  "

$ws.Range("D47").Value = "'" + "0.0898"
$ws.Range("E47").Value = "  -2.96%  "

$ws.Range("E48").Value = "  -4.38%  "

$ws.Range("E49").Value = "  -4.09%  "

$ws.Range("D50").Value = "'" + "16.49"
$ws.Range("E50").Value = "  -4.69%  "

$ws.Range("D51").Value = "1.696.96"
$ws.Range("E51").Value = "  -3.67%  "
